$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric values (e.g. "579.59") are not auto-converted to numbers
# by forcing the Text number format on target cells before assigning their values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.994.22'
$ws.Range("E2").Value = '  +4.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.254.20'
$ws.Range("E3").Value = '  +2.62%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.59'
$ws.Range("E5").Value = '  +3.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.04'
$ws.Range("E6").Value = '  +4.19%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").Value = '  -1.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.254.27'

$ws.Range("E10").Value = '  +4.37%  '

$ws.Range("E11").Value = '  +2.22%  '

$ws.Range("E12").Value = '  +4.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.822.06'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.136'
$ws.Range("E14").Value = '  +0.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.07'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.967.58'
$ws.Range("E16").Value = '  +4.33%  '

$ws.Range("E17").Value = '  +3.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.256.55'
$ws.Range("E18").Value = '  +2.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.86'
$ws.Range("E19").Value = '  +2.80%  '

$ws.Range("E20").Value = '  +2.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.81'
$ws.Range("E21").Value = '  +5.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.62'
$ws.Range("E22").Value = '  +6.17%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.42'
$ws.Range("E24").Value = '  +3.37%  '

$ws.Range("E25").Value = '  +2.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.402.97'
$ws.Range("E26").Value = '  +2.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000118'
$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("E28").Value = '  +2.43%  '

$ws.Range("E29").Value = '  +2.02%  '

$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("E31").Value = '  +4.25%  '

$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.62'
$ws.Range("E33").Value = '  +2.42%  '

$ws.Range("E34").Value = '  +0.10%  '

$ws.Range("E35").Value = '  +5.09%  '

$ws.Range("E36").Value = '  +2.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.76'
$ws.Range("E37").Value = '  +6.82%  '

$ws.Range("E38").Value = '  +4.60%  '

$ws.Range("E39").Value = '  +5.68%  '

$ws.Range("E40").Value = '  +9.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.12'
$ws.Range("E41").Value = '  +4.53%  '

$ws.Range("E42").Value = '  +2.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.48'
$ws.Range("E43").Value = '  +7.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.747.81'
$ws.Range("E44").Value = '  +5.58%  '

$ws.Range("E45").Value = '  +4.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '348.95'
$ws.Range("E46").Value = '  +5.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.13'
$ws.Range("E47").Value = '  +5.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.52'
$ws.Range("E48").Value = '  +4.32%  '

$ws.Range("E49").Value = '  +3.19%  '

$ws.Range("E50").Value = '  +3.43%  '

$ws.Range("E51").Value = '  +1.35%  '
